$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 4477.25
$ws.Range("I16").Value = 4477.25
$ws.Range("K16").Value = 4477.25
$ws.Range("M16").Value = -4247.25

$ws.Range("H135").Value = 2575.509
$ws.Range("I135").Value = 2329.8108
$ws.Range("J135").Value = 3080.5557
$ws.Range("K135").Value = 20968.2972
$ws.Range("L135").Value = 27725.0013
$ws.Range("M135").Value = -18433.2972
$ws.Range("N135").Value = -32795.0013

$ws.Range("H137").Value = 1612.96
$ws.Range("I137").Value = 1174.5555
$ws.Range("K137").Value = 3523.6665
$ws.Range("M137").Value = -973.6664999999998

$ws.Range("H138").Value = 2422.1797
$ws.Range("I138").Value = 1541.8846
$ws.Range("J138").Value = 2785.476
$ws.Range("K138").Value = 4625.6538
$ws.Range("L138").Value = 8356.428
$ws.Range("M138").Value = 514.3462
$ws.Range("N138").Value = -18636.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8113.7866
$ws.Range("I32").Value = 7190.0283
$ws.Range("K32").Value = 7190.0283
$ws.Range("M32").Value = -6903.0283

$ws.Range("H61").Value = 196703.98
$ws.Range("I61").Value = 5339.6333
$ws.Range("J61").Value = 457655.38
$ws.Range("K61").Value = 5339.6333
$ws.Range("L61").Value = 457655.38
$ws.Range("M61").Value = -5127.6333
$ws.Range("N61").Value = -458079.38

$ws.Range("H74").Value = 1960.3334
$ws.Range("I74").Value = 1576.8334
$ws.Range("K74").Value = 1576.8334
$ws.Range("M74").Value = -702.8334

$ws.Range("H77").Value = 1960.3334
$ws.Range("I77").Value = 1576.8334
$ws.Range("K77").Value = 7884.166999999999
$ws.Range("M77").Value = -3516.166999999999

$ws.Range("H107").Value = 50228
$ws.Range("J107").Value = 50228
$ws.Range("L107").Value = 50228
$ws.Range("N107").Value = -57908

$ws.Range("H122").Value = 714666
$ws.Range("I122").Value = 951951.0600000001
$ws.Range("J122").Value = 2811
$ws.Range("K122").Value = 2855853.18
$ws.Range("L122").Value = 8433
$ws.Range("M122").Value = -2853403.18
$ws.Range("N122").Value = -13333

$ws.Range("H132").Value = 2860680
$ws.Range("I132").Value = 2510.75
$ws.Range("K132").Value = 7532.25
$ws.Range("M132").Value = -5002.25

$ws.Range("H136").Value = 196703.98
$ws.Range("I136").Value = 5339.6333
$ws.Range("J136").Value = 457655.38
$ws.Range("K136").Value = 16018.8999
$ws.Range("L136").Value = 1372966.14
$ws.Range("M136").Value = -13468.8999
$ws.Range("N136").Value = -1378066.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 68000
$ws.Range("J122").Value = 68000
$ws.Range("L122").Value = 68000
$ws.Range("N122").Value = -77800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 90002
$ws.Range("J4").Value = 90002
$ws.Range("L4").Value = 90002
$ws.Range("N4").Value = -90226

$ws.Range("H22").Value = 481
$ws.Range("I22").Value = 481
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 481
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -131
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 7011.676
$ws.Range("I31").Value = 4704
$ws.Range("J31").Value = 7458.3228
$ws.Range("K31").Value = 4704
$ws.Range("L31").Value = 7458.3228
$ws.Range("M31").Value = -4409
$ws.Range("N31").Value = -8048.3228

$ws.Range("H34").Value = 7011.676
$ws.Range("I34").Value = 4704
$ws.Range("J34").Value = 7458.3228
$ws.Range("K34").Value = 4704
$ws.Range("L34").Value = 7458.3228
$ws.Range("M34").Value = -4502
$ws.Range("N34").Value = -7862.3228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 526371.8
$ws.Range("I8").Value = 526371.8
$ws.Range("K8").Value = 1579115.4
$ws.Range("M8").Value = -1578976.4

$ws.Range("H86").Value = 736
$ws.Range("I86").Value = 736
$ws.Range("K86").Value = 2208
$ws.Range("M86").Value = -1022

$ws.Range("H89").Value = 736
$ws.Range("I89").Value = 736
$ws.Range("K89").Value = 6624
$ws.Range("M89").Value = -696

$ws.Range("H92").Value = 827.8333
$ws.Range("I92").Value = 709.8333
$ws.Range("J92").Value = 842.5833
$ws.Range("K92").Value = 2129.4999
$ws.Range("L92").Value = 2527.7499
$ws.Range("M92").Value = -881.4998999999998
$ws.Range("N92").Value = -5023.7499

$ws.Range("H107").Value = 488.26315
$ws.Range("I107").Value = 426.2
$ws.Range("J107").Value = 721
$ws.Range("K107").Value = 1278.6
$ws.Range("L107").Value = 2163
$ws.Range("M107").Value = 641.4000000000001
$ws.Range("N107").Value = -6003

$ws.Range("H113").Value = 1935965.2
$ws.Range("I113").Value = 2381424.5
$ws.Range("J113").Value = 1000500.9
$ws.Range("K113").Value = 7144273.5
$ws.Range("L113").Value = 3001502.7
$ws.Range("M113").Value = -7142103.5
$ws.Range("N113").Value = -3005842.7

$ws.Range("H114").Value = 5532.8
$ws.Range("I114").Value = 273.5
$ws.Range("J114").Value = 9039
$ws.Range("K114").Value = 820.5
$ws.Range("L114").Value = 27117
$ws.Range("M114").Value = 2433.5
$ws.Range("N114").Value = -33625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H7").Value = 25000000
$ws.Range("I7").Value = 25000000
$ws.Range("K7").Value = 25000000
$ws.Range("M7").Value = -24999888

$ws.Range("H8").Value = 25000000
$ws.Range("I8").Value = 25000000
$ws.Range("K8").Value = 25000000
$ws.Range("M8").Value = -24999861

$ws.Range("H14").Value = 5502500
$ws.Range("I14").Value = 7333333.5
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 7333333.5
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -7333165.5
$ws.Range("N14").Value = -10336

$ws.Range("H113").Value = 40001390
$ws.Range("I113").Value = 47620084
$ws.Range("J113").Value = 3275
$ws.Range("K113").Value = 47620084
$ws.Range("L113").Value = 3275
$ws.Range("M113").Value = -47617914
$ws.Range("N113").Value = -7615

$ws.Range("H122").Value = 62640212
$ws.Range("I122").Value = 118315090
$ws.Range("J122").Value = 5975.375
$ws.Range("K122").Value = 354945270
$ws.Range("L122").Value = 17926.125
$ws.Range("M122").Value = -354942820
$ws.Range("N122").Value = -22826.125

$ws.Range("H123").Value = 28914.334
$ws.Range("J123").Value = 28914.334
$ws.Range("L123").Value = 28914.334
$ws.Range("N123").Value = -33814.334

$ws.Range("H132").Value = 3989.342
$ws.Range("I132").Value = 5081.3657
$ws.Range("J132").Value = 2710.1143
$ws.Range("K132").Value = 15244.0971
$ws.Range("L132").Value = 8130.342900000001
$ws.Range("M132").Value = -12714.0971
$ws.Range("N132").Value = -13190.3429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4530042.5
$ws.Range("I122").Value = 4768851.5
$ws.Range("K122").Value = 14306554.5
$ws.Range("M122").Value = -14304104.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 35667.668
$ws.Range("J2").Value = 35667.668
$ws.Range("L2").Value = 35667.668
$ws.Range("N2").Value = -35891.668

$ws.Range("H132").Value = 1880.909
$ws.Range("I132").Value = 1510.9584
$ws.Range("J132").Value = 2867.4443
$ws.Range("K132").Value = 4532.8752
$ws.Range("L132").Value = 8602.332900000001
$ws.Range("M132").Value = -2002.8752
$ws.Range("N132").Value = -13662.3329
